# Add "Sheet2" after the existing "Sheet1", populate it with the refreshed
# zip-code positivity-rate data (9/29 to 10/10), and update the selection /
# active-tab state to mirror the source commit.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "Sheet2"

# ---- Header row ----
$ws2.Range("B1").Value2 = "total"
$ws2.Range("C1").Value2 = "white"
$ws2.Range("D1").Value2 = "norm. diff."
$ws2.Range("E1").Value2 = "direction"

# ---- Data rows (zip, total, white, region, direction note) ----
$rows = @(
    @(60631, 4.7, 6.5, "west", $null),
    @(60656, 7.8, 10.4, "west", $null),
    @(60706, 5.5, 9.8, "west", $null),
    @(60634, 6.7, 6.3, "west", $null),
    @(60714, 4.1, 9, "north", $null),
    @(60053, 3.4, 5.1, "north", $null),
    @(60077, 0.9, 1.2, "north", $null),
    @(60076, 5.6, 5.7, "north", $null),
    @(60203, 3.9, 4, "north", $null),
    @(60712, 3.4, 1.8, "north", $null),
    @(60646, 3.4, 4.6, "north", "^ less urban"),
    @(60630, 6.3, 5.4, "my zip code", $null),
    @(60641, 4.8, 3.4, "south/east", "v more urban"),
    @(60645, 2.8, 4.8, "south/east", "heavily orthodox jewish"),
    @(60659, 6.3, 9.1, "south/east", $null),
    @(60625, 4.5, 4.1, "south/east", $null),
    @(60618, 3.9, 2.5, "south/east", $null)
)

$r = 2
foreach ($row in $rows) {
    $ws2.Cells.Item($r, 1).Value2 = $row[0]
    $ws2.Cells.Item($r, 2).Value2 = $row[1]
    $ws2.Cells.Item($r, 3).Value2 = $row[2]
    $ws2.Cells.Item($r, 5).Value2 = $row[3]
    if ($row[4]) {
        $ws2.Cells.Item($r, 6).Value2 = $row[4]
    }
    $r++
}

# D2 gets its own formula; D3:D18 becomes a shared formula, just like Sheet1.
$ws2.Range("D2").Formula = "=LN(C2/B2)"
$ws2.Range("D3:D18").Formula = "=LN(C3/B3)"

$ws2.Range("B19").Value2 = "data from 9/29 to 10/10"

# ---- Conditional formatting (same rules as Sheet1, different pri/order) ----
$null = $ws2.Range("D2:D18").FormatConditions.AddColorScale(3)
$null = $ws2.Range("B2:C18").FormatConditions.AddColorScale(2)

$cfD = $ws2.Range("D2:D18").FormatConditions(1)
$cfD.ColorScaleCriteria(1).Type = 1
$cfD.ColorScaleCriteria(1).FormatColor.Color = 8109667
$cfD.ColorScaleCriteria(2).Type = 0
$cfD.ColorScaleCriteria(2).Value = 0
$cfD.ColorScaleCriteria(2).FormatColor.Color = 16776444
$cfD.ColorScaleCriteria(3).Type = 2
$cfD.ColorScaleCriteria(3).FormatColor.Color = 7039480

$cfB = $ws2.Range("B2:C18").FormatConditions(1)
$cfB.ColorScaleCriteria(1).Type = 1
$cfB.ColorScaleCriteria(1).FormatColor.Color = 16776444
$cfB.ColorScaleCriteria(2).Type = 2
$cfB.ColorScaleCriteria(2).FormatColor.Color = 7039480

$cfB.Priority = 1
$cfD.Priority = 2

# ---- Selections / active tab ----
$sheet1.Activate()
$sheet1.Range("D29").Select() | Out-Null

$ws2.Activate()
$ws2.Range("E10").Select() | Out-Null
